$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '27.435.77'
$ws.Range('E2').Value = '  +0.12%  '
$ws.Range('D3').Value = '1.636.20'
$ws.Range('E3').Value = '  -0.91%  '
$ws.Range('E4').Value = '  +0.02%  '
$ws.Range('D5').Value = "'212.37"
$ws.Range('E5').Value = '  -0.71%  '
$ws.Range('D6').Value = "'0.533"
$ws.Range('E6').Value = '  +4.82%  '
$ws.Range('E7').Value = '  +0.02%  '
$ws.Range('D8').Value = "'22.96"
$ws.Range('E8').Value = '  -5.27%  '
$ws.Range('E9').Value = '  -2.34%  '
$ws.Range('E10').Value = '  -0.92%  '
$ws.Range('D11').Value = "'0.0886"
$ws.Range('E11').Value = '  +1.16%  '
$ws.Range('D12').Value = '1.867.54'
$ws.Range('E12').Value = '  -0.91%  '
$ws.Range('D13').Value = '1.635.15'
$ws.Range('E13').Value = '  -1.03%  '
$ws.Range('D14').Value = "'0.580"
$ws.Range('E14').Value = '  +3.04%  '
$ws.Range('E15').Value = '  -1.75%  '
$ws.Range('D16').Value = "'64.21"
$ws.Range('E16').Value = '  -2.29%  '
$ws.Range('D17').Value = '27.434.52'
$ws.Range('E17').Value = '  +0.13%  '
$ws.Range('D18').Value = "'228.96"
$ws.Range('E18').Value = '  -2.55%  '
$ws.Range('E19').Value = '  -0.41%  '
$ws.Range('D20').Value = "'7.61"
$ws.Range('E20').Value = '  +1.16%  '
$ws.Range('E21').Value = '  -0.11%  '
$ws.Range('E22').Value = '  -2.30%  '
$ws.Range('E23').Value = '  +4.90%  '
$ws.Range('E24').Value = '  -3.43%  '
$ws.Range('D25').Value = "'149.51"
$ws.Range('E25').Value = '  +2.37%  '
$ws.Range('D26').Value = "'6.99"
$ws.Range('E26').Value = '  -2.71%  '
$ws.Range('E27').Value = '  +1.81%  '
$ws.Range('E28').Value = '  +0.00%  '
$ws.Range('D29').Value = "'15.54"
$ws.Range('E29').Value = '  -3.28%  '
$ws.Range('E30').Value = '  -0.77%  '
$ws.Range('E31').Value = '  -1.96%  '
$ws.Range('E32').Value = '  -0.42%  '
$ws.Range('E33').Value = '  +2.82%  '
$ws.Range('D34').Value = '1.420.95'
$ws.Range('E34').Value = '  -3.05%  '
$ws.Range('E35').Value = '  +2.49%  '
$ws.Range('E36').Value = '  -1.92%  '
$ws.Range('D37').Value = "'0.571"
$ws.Range('E37').Value = '  +0.09%  '
$ws.Range('E38').Value = '  -1.19%  '
$ws.Range('E39').Value = '  -4.16%  '
$ws.Range('D40').Value = "'0.891"
$ws.Range('E40').Value = '  +13.25%  '
$ws.Range('E41').Value = '  -1.56%  '
$ws.Range('E42').Value = '  +0.04%  '
$ws.Range('E43').Value = '  +1.46%  '
$ws.Range('D45').Value = "'64.88"
$ws.Range('E45').Value = '  -0.48%  '
$ws.Range('D46').Value = '1.777.50'
$ws.Range('E46').Value = '  -0.78%  '
$ws.Range('E47').Value = '  -3.32%  '
$ws.Range('D48').Value = "'85.83"
$ws.Range('E48').Value = '  -2.82%  '
$ws.Range('D49').Value = '0.0₆0107'
$ws.Range('E49').Value = '  +0.10%  '
$ws.Range('D50').Value = "'0.0989"
$ws.Range('E50').Value = '  -1.94%  '
$ws.Range('D51').Value = "'7.70"
$ws.Range('E51').Value = '  -1.16%  '
